$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 67
$templateRow = 66

# Copy the formatting from the row above (so the new date cell keeps the
# same date number-format / font / border style as the existing rows)
$ws.Cells.Item($templateRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 44775
$ws.Cells.Item($newRow, 2).Value = 751.37
$ws.Cells.Item($newRow, 3).Value = 11889.68
$ws.Cells.Item($newRow, 4).Value = 1853.39
